$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.428.48'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.24%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.864.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.76%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.22%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.7060'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.87%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3135'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.88%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07843'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.41'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.79%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08009'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -4.03%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.888.22'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.192'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.61%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6998'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.49%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.452'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.24%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.482.43'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.38%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008339'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.59%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '252.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.66%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.138.67'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.46%  '

$ws.Range("E21").Value = '  -1.40%  '

$ws.Range("E22").Value = '  -0.01%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.598'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.81%  '

$ws.Range("E24").Value = '  -0.06%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1554'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.26%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.999'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '160.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.50%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.61%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.499'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.64%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.317'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.272'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.205'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.02%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05299'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.05%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.877'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.55%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7507'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.82%  '

$ws.Range("E36").Value = '  -2.06%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.715'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01875'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.83%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.258.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.739'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8984'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.52%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '108.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.05%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.964'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -8.50%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '71.22'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.55%  '

$ws.Range("E45").Value = '  -0.01%  '

$ws.Range("E46").Value = '  -0.56%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.037.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.16%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.5191'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.61%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.784'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.52%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.515'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4302'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.87%  '
